# Horizontally center the page numbers in the document footer.
$d = $word.ActiveDocument

# The document's primary footer contains the page-number field in its
# first paragraph. Center that paragraph (wdAlignParagraphCenter = 1).
$footer = $d.Sections(1).Footers(1)
$p = $footer.Range.Paragraphs(1)
$p.Alignment = 1
